$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Generate Report for Handback
#
# The localization-status report is regenerated after a handback completes:
#   - Status flips from "Ready for handoff" to "Handed back: in sync with en-US"
#     (this text is a shared string, so it ripples through every sheet that
#     shows a Status/zh-cn/de-de column).
#   - Each language sheet (zh-cn, de-de) gets its "Latest Target File" (I) and
#     "Latest Handback File" (J) columns populated with the source markdown
#     file name and the generated xlf file name, each becoming a hyperlink.
#   - "Latest Handback DateTime" (K) moves off the zero date to the real
#     handback timestamp.
# ---------------------------------------------------------------------------

$mdFile1 = "2157b7b0-8046-4c28-a10f-5ff9a7c16c57.md"
$mdFile2 = "228a2111-de56-4b94-b7db-5cad7f67987a.md"

$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/oltest/blob/343e9a8000f68139d9377ea5b6e3b59f3b762263/e2e/2157b7b0-8046-4c28-a10f-5ff9a7c16c57.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/oltest/blob/343e9a8000f68139d9377ea5b6e3b59f3b762263/e2e/228a2111-de56-4b94-b7db-5cad7f67987a.md"

$xlfZhCn1 = "2157b7b0-8046-4c28-a10f-5ff9a7c16c57.8b6e386a9bf098c08e8093e447e2146ccadef719.zh-cn.xlf"
$xlfZhCn2 = "228a2111-de56-4b94-b7db-5cad7f67987a.0abc0f6c5aadc0ffe4e2ad1336e44de04523adfd.zh-cn.xlf"
$xlfDeDe1 = "2157b7b0-8046-4c28-a10f-5ff9a7c16c57.8b6e386a9bf098c08e8093e447e2146ccadef719.de-de.xlf"
$xlfDeDe2 = "228a2111-de56-4b94-b7db-5cad7f67987a.0abc0f6c5aadc0ffe4e2ad1336e44de04523adfd.de-de.xlf"

$handedBackStatus = "Handed back: in sync with en-US"
$handbackTimeZhCn = "2016-08-12 05:04:52"
$handbackTimeDeDe = "2016-08-12 05:05:00"

$newColWidth = 29.166666666666668   # renders as 29.9777047293527 -> nearest grid value (40 for capped cols handled separately)
$wideColWidth = 39.166666666666664  # renders as 40 (the workbook's long-text cap)

# ---------------------------------------------------------------------------
# Overview sheet: the Status text ("Ready for handoff" -> handed-back text)
# is shared, so editing it on the language sheets below already updates the
# strings everywhere they're referenced. We only need to touch the Overview
# column widths that widen to fit the longer status text.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $newColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColWidth

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = $handedBackStatus
$wsZhCn.Range("C3").Value = $handedBackStatus

$wsZhCn.Range("I2").Value = $mdFile1
$wsZhCn.Range("J2").Value = $xlfZhCn1
$wsZhCn.Range("K2").Value = $handbackTimeZhCn

$wsZhCn.Range("I3").Value = $mdFile2
$wsZhCn.Range("J3").Value = $xlfZhCn2
$wsZhCn.Range("K3").Value = $handbackTimeZhCn

# Recreate the hyperlinks in display order (A2, I2, A3, I3) so the new links
# land on I2/I3 and the relationship ids line up the same way Excel assigns
# them when hyperlinks are (re)written in sheet order.
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $mdUrl1, "", "", $mdFile1)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdUrl1, "", "", $mdFile1)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $mdUrl2, "", "", $mdFile2)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $mdUrl2, "", "", $mdFile2)

$wsZhCn.Columns.Item(3).ColumnWidth = $newColWidth
$wsZhCn.Columns.Item(9).ColumnWidth = $wideColWidth
$wsZhCn.Columns.Item(10).ColumnWidth = $wideColWidth

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = $handedBackStatus
$wsDeDe.Range("C3").Value = $handedBackStatus

$wsDeDe.Range("I2").Value = $mdFile1
$wsDeDe.Range("J2").Value = $xlfDeDe1
$wsDeDe.Range("K2").Value = $handbackTimeDeDe

$wsDeDe.Range("I3").Value = $mdFile2
$wsDeDe.Range("J3").Value = $xlfDeDe2
$wsDeDe.Range("K3").Value = $handbackTimeDeDe

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $mdUrl1, "", "", $mdFile1)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdUrl1, "", "", $mdFile1)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $mdUrl2, "", "", $mdFile2)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $mdUrl2, "", "", $mdFile2)

$wsDeDe.Columns.Item(3).ColumnWidth = $newColWidth
$wsDeDe.Columns.Item(9).ColumnWidth = $wideColWidth
$wsDeDe.Columns.Item(10).ColumnWidth = $wideColWidth

Write-Host "Handback report regenerated"
